$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.986.37"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.56%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.501.18"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.44%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.84"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.40%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.62"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.93%  "

# Row 7
$ws.Range("E7").Value = "  +0.03%  "

# Row 8
$ws.Range("E8").Value = "  -0.24%  "

# Row 9
$ws.Range("E9").Value = "  +3.58%  "

# Row 10
$ws.Range("E10").Value = "  -1.37%  "

# Row 11
$ws.Range("E11").Value = "  -0.88%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.105.68"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.40%  "

# Row 13
$ws.Range("E13").Value = "  -0.18%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "29.20"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +3.81%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "66.985.27"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.54%  "

# Row 16
$ws.Range("E16").Value = "  +0.46%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.538.87"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.66%  "

# Row 18
$ws.Range("E18").Value = "  -0.29%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.28"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.78%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "396.04"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.93%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.99"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.05%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.46"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.63%  "

# Row 23
$ws.Range("E23").Value = "  +0.08%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.536"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.30%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000122"
$ws.Range("D25").ClearFormats()

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.26"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.65%  "

# Row 27
$ws.Range("E27").Value = "  +0.30%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.997"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.23%  "

# Row 29
$ws.Range("E29").Value = "  -1.93%  "

# Row 30
$ws.Range("E30").Value = "  -2.19%  "

# Row 31
$ws.Range("E31").Value = "  -0.11%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "23.76"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.94%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.38"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.36%  "

# Row 34
$ws.Range("E34").Value = "  +1.15%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "162.79"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.83%  "

# Row 36
$ws.Range("E36").Value = "  -2.59%  "

# Row 37
$ws.Range("E37").Value = "  -0.24%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.92"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.00%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.66"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.40%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0741"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.58%  "

# Row 41
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.839.29"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.66%  "

# Row 42
$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "27.16"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.59%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.30"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.08%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.79"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.91%  "

# Row 45
$ws.Range("E45").Value = "  +2.42%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0302"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.31%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "338.69"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.47%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "34.64"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.08%  "

# Row 49
$ws.Range("E49").Value = "  -0.99%  "

# Row 50
$ws.Range("E50").Value = "  -4.53%  "

# Row 51
$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.43"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.34%  "
